$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rooms")

$ws.Range("A1").Value = "Room#"
$ws.Range("B1").Value = "Room Type"

$ws.Range("K14").Select()
